$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IR LED row (B4): add hyperlink (text already shows the URL), style it like the
# other "Link" cells (B3/B6).
$irUrl = "https://www.lcsc.com/product-detail/span-style-background-color-ff0-Infrared-span-IR-LEDs_EKINGLUX-E6QYDD1204-IRA940nm_C396649.html"
$ws.Hyperlinks.Add($ws.Range("B4"), $irUrl)
$ws.Range("B4").Style = "Link"

# 590 Ohm Resistor row (B8): replace the old (now obsolete) LCSC link with the
# new part, and turn it into a real hyperlink too.
$resistorUrl = "https://www.lcsc.com/product-detail/Chip-span-style-background-color-ff0-Resistor-span-Surface-Mount_Sunway-SC0603J5600F2BNRH_C5140946.html"
$ws.Range("B8").Value = $resistorUrl
$ws.Hyperlinks.Add($ws.Range("B8"), $resistorUrl)
$ws.Range("B8").Style = "Link"

# Move the active selection from B3 to B8.
[void]$ws.Range("B8").Select()
